$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 520 ("「開かれたターと結ばれたター」"), shifting all rows below up by one.
$ws.Rows.Item(520).Delete()
